$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change PO customer_code / customer_name headers to supplier_code / supplier_name
$ws.Range("C1").Value = "Supplier Code"
$ws.Range("D1").Value = "Supplier Name"

# Update selection to D1 (matches the new state captured in the sheet view)
$ws.Range("D1").Select()
